$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1736781
$ws.Range("I55").Value = 1308.625
$ws.Range("K55").Value = 1308.625
$ws.Range("M55").Value = -1094.625
$ws.Range("H62").Value = 4934.8237
$ws.Range("I62").Value = 4854.4443
$ws.Range("J62").Value = 5025.25
$ws.Range("K62").Value = 4854.4443
$ws.Range("L62").Value = 5025.25
$ws.Range("M62").Value = -4230.4443
$ws.Range("N62").Value = -6273.25
$ws.Range("H64").Value = 2692.93
$ws.Range("I64").Value = 2552.5527
$ws.Range("J64").Value = 2973.6843
$ws.Range("K64").Value = 2552.5527
$ws.Range("L64").Value = 2973.6843
$ws.Range("M64").Value = -2304.5527
$ws.Range("N64").Value = -3469.6843
$ws.Range("H65").Value = 4934.8237
$ws.Range("I65").Value = 4854.4443
$ws.Range("J65").Value = 5025.25
$ws.Range("K65").Value = 24272.2215
$ws.Range("L65").Value = 25126.25
$ws.Range("M65").Value = -21152.2215
$ws.Range("N65").Value = -31366.25
$ws.Range("H67").Value = 2692.93
$ws.Range("I67").Value = 2552.5527
$ws.Range("J67").Value = 2973.6843
$ws.Range("K67").Value = 2552.5527
$ws.Range("L67").Value = 2973.6843
$ws.Range("M67").Value = -1694.5527
$ws.Range("N67").Value = -4689.6843

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1999.7693
$ws.Range("I3").Value = 274.25
$ws.Range("K3").Value = 274.25
$ws.Range("M3").Value = -159.25
$ws.Range("H4").Value = 79
$ws.Range("I4").Value = 79
$ws.Range("K4").Value = 79
$ws.Range("M4").Value = 37
$ws.Range("H5").Value = 182.66667
$ws.Range("J5").Value = 224
$ws.Range("L5").Value = 224
$ws.Range("N5").Value = -448
$ws.Range("H45").Value = 43032.168
$ws.Range("I45").Value = 101207.4
$ws.Range("J45").Value = 1478.4286
$ws.Range("K45").Value = 101207.4
$ws.Range("L45").Value = 1478.4286
$ws.Range("M45").Value = -100830.4
$ws.Range("N45").Value = -2232.4286
$ws.Range("H61").Value = 1442.3103
$ws.Range("I61").Value = 1418.7778
$ws.Range("J61").Value = 1480.8182
$ws.Range("K61").Value = 1418.7778
$ws.Range("L61").Value = 1480.8182
$ws.Range("M61").Value = -1206.7778
$ws.Range("N61").Value = -1904.8182
$ws.Range("H74").Value = 20001352
$ws.Range("I74").Value = 22728678
$ws.Range("J74").Value = 966.3333
$ws.Range("K74").Value = 22728678
$ws.Range("L74").Value = 966.3333
$ws.Range("M74").Value = -22727804
$ws.Range("N74").Value = -2714.3333
$ws.Range("H77").Value = 20001352
$ws.Range("I77").Value = 22728678
$ws.Range("J77").Value = 966.3333
$ws.Range("K77").Value = 113643390
$ws.Range("L77").Value = 4831.6665
$ws.Range("M77").Value = -113639022
$ws.Range("N77").Value = -13567.6665
$ws.Range("H136").Value = 1442.3103
$ws.Range("I136").Value = 1418.7778
$ws.Range("J136").Value = 1480.8182
$ws.Range("K136").Value = 4256.3334
$ws.Range("L136").Value = 4442.4546
$ws.Range("M136").Value = -1706.3334
$ws.Range("N136").Value = -9542.454600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 182.66667
$ws.Range("J4").Value = 224
$ws.Range("L4").Value = 224
$ws.Range("N4").Value = -454
$ws.Range("H22").Value = 3378643.2
$ws.Range("I22").Value = 3378643.2
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3378643.2
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -3378470.2
$ws.Range("N22").ClearContents()
$ws.Range("H134").Value = 3836096.5
$ws.Range("I134").Value = 1223
$ws.Range("J134").Value = 13902639
$ws.Range("K134").Value = 3669
$ws.Range("L134").Value = 41707917
$ws.Range("M134").Value = -1134
$ws.Range("N134").Value = -41712987

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1333.3334
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -887
$ws.Range("N3").Value = -2226
$ws.Range("H22").Value = 9383.416999999999
$ws.Range("I22").Value = 243.28572
$ws.Range("K22").Value = 243.28572
$ws.Range("M22").Value = 106.71428
$ws.Range("H122").Value = 41667452
$ws.Range("I122").Value = 41667452
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 125002356
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -124999906
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 16668737
$ws.Range("I132").Value = 1830
$ws.Range("J132").Value = 47621564
$ws.Range("K132").Value = 5490
$ws.Range("L132").Value = 142864692
$ws.Range("M132").Value = -2960
$ws.Range("N132").Value = -142869752
$ws.Range("H134").Value = 2656.6667
$ws.Range("I134").Value = 1875.5555
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 5626.666499999999
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -3091.666499999999
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 37638224
$ws.Range("I5").Value = 41667004
$ws.Range("J5").Value = 33340858
$ws.Range("K5").Value = 125001012
$ws.Range("L5").Value = 100022574
$ws.Range("M5").Value = -125000900
$ws.Range("N5").Value = -100022798
$ws.Range("H20").Value = 625875
$ws.Range("H40").Value = 578.41174
$ws.Range("I40").Value = 158.25
$ws.Range("J40").Value = 707.6923
$ws.Range("K40").Value = 633
$ws.Range("L40").Value = 2830.7692
$ws.Range("M40").Value = -564
$ws.Range("N40").Value = -2968.7692
$ws.Range("H135").Value = 37638224
$ws.Range("I135").Value = 41667004
$ws.Range("J135").Value = 33340858
$ws.Range("K135").Value = 375003036
$ws.Range("L135").Value = 300067722
$ws.Range("M135").Value = -375000501
$ws.Range("N135").Value = -300072792

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 59.0625
$ws.Range("I2").Value = 11.666667
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 11.666667
$ws.Range("L2").Value = 120
$ws.Range("M2").Value = 101.333333
$ws.Range("N2").Value = -346
$ws.Range("H4").Value = 800
$ws.Range("J4").Value = 1250
$ws.Range("L4").Value = 1250
$ws.Range("N4").Value = -1474
$ws.Range("H126").Value = 2077.6316
$ws.Range("I126").Value = 1541.909
$ws.Range("J126").Value = 2814.25
$ws.Range("K126").Value = 4625.727000000001
$ws.Range("L126").Value = 8442.75
$ws.Range("M126").Value = -2155.727000000001
$ws.Range("N126").Value = -13382.75
$ws.Range("H132").Value = 5393.2905
$ws.Range("I132").Value = 1409.7
$ws.Range("K132").Value = 4229.1
$ws.Range("M132").Value = -1699.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2473.9524
$ws.Range("I22").Value = 962.625
$ws.Range("J22").Value = 3404
$ws.Range("K22").Value = 962.625
$ws.Range("L22").Value = 3404
$ws.Range("M22").Value = -667.625
$ws.Range("N22").Value = -3994
$ws.Range("H27").Value = 2473.9524
$ws.Range("I27").Value = 962.625
$ws.Range("J27").Value = 3404
$ws.Range("K27").Value = 962.625
$ws.Range("L27").Value = 3404
$ws.Range("M27").Value = -855.625
$ws.Range("N27").Value = -3618
$ws.Range("H55").Value = 11493.444
$ws.Range("I55").Value = 20260.2
$ws.Range("J55").Value = 535
$ws.Range("K55").Value = 20260.2
$ws.Range("L55").Value = 535
$ws.Range("M55").Value = -20087.2
$ws.Range("N55").Value = -881
$ws.Range("H122").Value = 6912.25
$ws.Range("I122").Value = 9421
$ws.Range("K122").Value = 28263
$ws.Range("M122").Value = -25813
$ws.Range("H132").Value = 37043880
$ws.Range("I132").Value = 66669780
$ws.Range("J132").Value = 11499
$ws.Range("K132").Value = 200009340
$ws.Range("L132").Value = 34497
$ws.Range("M132").Value = -200006810
$ws.Range("N132").Value = -39557

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4057
$ws.Range("I62").Value = 3937.5
$ws.Range("J62").Value = 4216.3335
$ws.Range("K62").Value = 3937.5
$ws.Range("L62").Value = 4216.3335
$ws.Range("M62").Value = -3313.5
$ws.Range("N62").Value = -5464.3335
$ws.Range("H65").Value = 4057
$ws.Range("I65").Value = 3937.5
$ws.Range("J65").Value = 4216.3335
$ws.Range("K65").Value = 19687.5
$ws.Range("L65").Value = 21081.6675
$ws.Range("M65").Value = -16567.5
$ws.Range("N65").Value = -27321.6675
